$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.090.06"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5145"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3750"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.56%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.18"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9046"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07645"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.890.86"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.274"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008474"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.44"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.111.37"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.069"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.135.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.66%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.402"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.05"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.782"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.11"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.58"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.963"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.837"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09181"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05092"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7825"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.288"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.630"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02001"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5583"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.989"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.633"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.80"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4801"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.25"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.596"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.09"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.87%  "
